$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows right after the header row, pushing all existing task
# rows down by 4 (row 2 -> row 6, ... row 23 -> row 27).
$ws.Rows("2:5").Insert()

# "di choi" task (now row 24): was logged as overdue (start 04-19, end 04-21,
# "Chưa hoàn thành"). Fix the overdue-count logic bug by correcting the
# dates/status to on-time completion.
$ws.Range("C24:D24").NumberFormat = "@"
$ws.Range("C24").Value = "2021-04-13"
$ws.Range("D24").Value = "2021-04-16"
$ws.Range("C24:D24").ClearFormats()
$ws.Range("E24").Value = "Hoàn thành"

# "zcx" task (now row 23): status corrected from "Chưa hoàn thành" to
# "Hoàn thành".
$ws.Range("E23").Value = "Hoàn thành"

# "hung" task (now row 27): end_date and status corrected.
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2021-04-10"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "Hoàn thành"
